$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update names and values (A2:B11) per the target diff
$ws.Range("A2").Value = "CARRION LAZARO MICHAEL LUIS"
$ws.Range("B2").Value = 138

$ws.Range("A3").Value = "ARRUNATEGUI ESPINOZA JOVANNY"
$ws.Range("B3").Value = 135

$ws.Range("A4").Value = "PAZ ANASTACIO JUANITA ROSA"
$ws.Range("B4").Value = 133

$ws.Range("A5").Value = "NIMA CARMEN KAREN DEL MILAGRO"
$ws.Range("B5").Value = 120

$ws.Range("A6").Value = "ALZAMORA CHERRES SIRLEY YASMIN"
$ws.Range("B6").Value = 114

$ws.Range("A7").Value = "ESPINOZA VALDIVIEZO JUNIOR RICARDO"
$ws.Range("B7").Value = 110

$ws.Range("A8").Value = "PULACHE LAZO VILMA YOHANA"
$ws.Range("B8").Value = 100

$ws.Range("A9").Value = "NAVARRO JUAREZ LIDIA"
$ws.Range("B9").Value = 97

$ws.Range("A10").Value = "DOMINGUEZ CUEVA MERLING DEL JESUS YOLINDA"
$ws.Range("B10").Value = 95

$ws.Range("A11").Value = "LILIAN ROXANA VEGA GARCÍA"
$ws.Range("B11").Value = 92
